$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 331
$ws.Range("I33").Value = 331
$ws.Range("K33").Value = 331
$ws.Range("M33").Value = -102

$ws.Range("H137").Value = 9804949
$ws.Range("I137").Value = 12821335
$ws.Range("J137").Value = 1696.0834
$ws.Range("K137").Value = 38464005
$ws.Range("L137").Value = 5088.2502
$ws.Range("M137").Value = -38461455
$ws.Range("N137").Value = -10188.2502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 734
$ws.Range("I4").Value = 600
$ws.Range("J4").Value = 1002
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 1002
$ws.Range("M4").Value = -484
$ws.Range("N4").Value = -1234

$ws.Range("H45").Value = 56774.832
$ws.Range("I45").Value = 125838.875
$ws.Range("J45").Value = 1523.6
$ws.Range("K45").Value = 125838.875
$ws.Range("L45").Value = 1523.6
$ws.Range("M45").Value = -125461.875
$ws.Range("N45").Value = -2277.6

$ws.Range("H74").Value = 938.9815
$ws.Range("I74").Value = 909.53656
$ws.Range("J74").Value = 1031.8462
$ws.Range("K74").Value = 909.53656
$ws.Range("L74").Value = 1031.8462
$ws.Range("M74").Value = -35.53656000000001
$ws.Range("N74").Value = -2779.8462

$ws.Range("H77").Value = 938.9815
$ws.Range("I77").Value = 909.53656
$ws.Range("J77").Value = 1031.8462
$ws.Range("K77").Value = 4547.6828
$ws.Range("L77").Value = 5159.231
$ws.Range("M77").Value = -179.6828000000005
$ws.Range("N77").Value = -13895.231

$ws.Range("H88").Value = 2504.96
$ws.Range("I88").Value = 2690.3845
$ws.Range("J88").Value = 2304.0833
$ws.Range("K88").Value = 2690.3845
$ws.Range("L88").Value = 2304.0833
$ws.Range("M88").Value = -2284.3845
$ws.Range("N88").Value = -3116.0833

$ws.Range("H91").Value = 2504.96
$ws.Range("I91").Value = 2690.3845
$ws.Range("J91").Value = 2304.0833
$ws.Range("K91").Value = 2690.3845
$ws.Range("L91").Value = 2304.0833
$ws.Range("M91").Value = -1286.3845
$ws.Range("N91").Value = -5112.0833

$ws.Range("H107").Value = 29038.666
$ws.Range("J107").Value = 29038.666
$ws.Range("L107").Value = 29038.666
$ws.Range("N107").Value = -36718.666

$ws.Range("H132").Value = 866198.44
$ws.Range("I132").Value = 912.625
$ws.Range("J132").Value = 4904199
$ws.Range("K132").Value = 2737.875
$ws.Range("L132").Value = 14712597
$ws.Range("M132").Value = -207.875
$ws.Range("N132").Value = -14717657

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1663200.4
$ws.Range("I86").Value = 1926.5714
$ws.Range("J86").Value = 3324474.2
$ws.Range("K86").Value = 1926.5714
$ws.Range("L86").Value = 3324474.2
$ws.Range("M86").Value = -803.5714
$ws.Range("N86").Value = -3326720.2

$ws.Range("H89").Value = 1663200.4
$ws.Range("I89").Value = 1926.5714
$ws.Range("J89").Value = 3324474.2
$ws.Range("K89").Value = 9632.857
$ws.Range("L89").Value = 16622371
$ws.Range("M89").Value = -4016.857
$ws.Range("N89").Value = -16633603

$ws.Range("H134").Value = 3067.9807
$ws.Range("I134").Value = 1045.8914
$ws.Range("J134").Value = 18570.666
$ws.Range("K134").Value = 3137.6742
$ws.Range("L134").Value = 55711.99800000001
$ws.Range("M134").Value = -602.6741999999999
$ws.Range("N134").Value = -60781.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()

$ws.Range("H31").Value = 2421240.8
$ws.Range("I31").Value = 6184959.5
$ws.Range("J31").Value = 1707.1428
$ws.Range("K31").Value = 6184959.5
$ws.Range("L31").Value = 1707.1428
$ws.Range("M31").Value = -6184664.5
$ws.Range("N31").Value = -2297.1428

$ws.Range("H34").Value = 2421240.8
$ws.Range("I34").Value = 6184959.5
$ws.Range("J34").Value = 1707.1428
$ws.Range("K34").Value = 6184959.5
$ws.Range("L34").Value = 1707.1428
$ws.Range("M34").Value = -6184757.5
$ws.Range("N34").Value = -2111.1428

$ws.Range("H62").Value = 3966.4285
$ws.Range("I62").Value = 2994.75
$ws.Range("J62").Value = 5262
$ws.Range("K62").Value = 2994.75
$ws.Range("L62").Value = 5262
$ws.Range("M62").Value = -2370.75
$ws.Range("N62").Value = -6510

$ws.Range("H65").Value = 3966.4285
$ws.Range("I65").Value = 2994.75
$ws.Range("J65").Value = 5262
$ws.Range("K65").Value = 14973.75
$ws.Range("L65").Value = 26310
$ws.Range("M65").Value = -11853.75
$ws.Range("N65").Value = -32550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 985
$ws.Range("I19").Value = 985
$ws.Range("K19").Value = 985
$ws.Range("M19").Value = -697

$ws.Range("H102").Value = 2662.5557
$ws.Range("I102").Value = 1319
$ws.Range("J102").Value = 3517.5454
$ws.Range("K102").Value = 1319
$ws.Range("L102").Value = 3517.5454
$ws.Range("M102").Value = 303
$ws.Range("N102").Value = -6761.5454

$ws.Range("H122").Value = 26322232
$ws.Range("I122").Value = 38470260
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 115410780
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -115408330
$ws.Range("N122").Value = -9398.5

$ws.Range("H132").Value = 4713.186
$ws.Range("I132").Value = 2132.4243
$ws.Range("J132").Value = 13229.7
$ws.Range("K132").Value = 6397.2729
$ws.Range("L132").Value = 39689.10000000001
$ws.Range("M132").Value = -3867.2729
$ws.Range("N132").Value = -44749.10000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 845916
$ws.Range("I22").Value = 1810042.6
$ws.Range("J22").Value = 2305.25
$ws.Range("K22").Value = 1810042.6
$ws.Range("L22").Value = 2305.25
$ws.Range("M22").Value = -1809747.6
$ws.Range("N22").Value = -2895.25

$ws.Range("H27").Value = 845916
$ws.Range("I27").Value = 1810042.6
$ws.Range("J27").Value = 2305.25
$ws.Range("K27").Value = 1810042.6
$ws.Range("L27").Value = 2305.25
$ws.Range("M27").Value = -1809935.6
$ws.Range("N27").Value = -2519.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H123").Value = 20333.334
$ws.Range("J123").Value = 20333.334
$ws.Range("L123").Value = 20333.334
$ws.Range("N123").Value = -30133.334
